# repull data, push all data, mean calculation
# Update the dSF column (F) with freshly pulled values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -2
    4  = -1
    5  = -3
    7  = 1
    9  = -1
    10 = -1
    11 = -3
    13 = -2
    15 = -4
    17 = -3
    19 = -5
    20 = -7
    21 = -3
    22 = -5
    23 = -1
    24 = -2
    25 = -3
    26 = 2
    27 = -3
    28 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
